$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the robotic S1 prep kit value (shared across G2:G27) to add the
#    "L" suffix: NEBNextPoly(A)E7490 -> NEBNextPoly(A)E7490L
$ws.Range("G2:G27").Value = "NEBNextPoly(A)E7490L"

# 2. Column G ("roboticS1Prep") now holds longer text, so widen it; keep the
#    other columns at their original width. (29 is the closest ColumnWidth
#    the host's pixel-quantized engine can render as ~29.89 XML width.)
$ws.Columns.Item(7).ColumnWidth = 29

# 3. Column I (the accuracy-check flag) switches from a literal boolean to an
#    explicit =FALSE() formula for every data row.
for ($i = 2; $i -le 27; $i++) {
    $ws.Range("I" + $i).Formula = "=FALSE()"
}

# 4. Move the active selection from the old I2:I27 focus to the newly
#    widened G2:G27 column that was just updated.
$ws.Range("G2:G27").Select()
